$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Documentation / units corrections: strip the decorative parentheses
#     (and stray quotes) that wrapped these unit labels in the shared
#     strings table. ---
$ws.Range("C6").Value  = "x10^4 \mu m"
$ws.Range("C7").Value  = "x10^4 \mu m^2"
$ws.Range("C8").Value  = "x10^4 \mu m^3"
$ws.Range("C15").Value = "\text{degrees}"
$ws.Range("C16").Value = "\text{degrees}"

# --- Column width touch-up (minor manual resize of the first three
#     columns that accompanied the text correction). ---
$ws.Columns.Item(1).ColumnWidth = 31.1666666666667
$ws.Columns.Item(2).ColumnWidth = 9.16666666666667
$ws.Columns.Item(3).ColumnWidth = 14.8333333333333

# --- Move the active selection to C16 ---
$ws.Range("C16").Select()
